# Add "User info" rows (street/neighborhood/city + unrequired key) to the
# language-values table, mirroring rows 3-4 of the existing sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - street name
$ws.Range("B5").Value = "system.common.streetName"
$ws.Range("D5").Value = "Street name"
$ws.Range("E5").Value = "Улица/блок"
$ws.Range("F5").Value = "Улица"

# Row 6 - neighborhood
$ws.Range("B6").Value = "system.common.Neighborhood"
$ws.Range("D6").Value = "Neighborhood"
$ws.Range("E6").Value = "Квартал"
$ws.Range("F6").Value = "Микрорайон"

# Row 7 - city
$ws.Range("B7").Value = "system.common.City"
$ws.Range("D7").Value = "City"
$ws.Range("E7").Value = "Град"
$ws.Range("F7").Value = "Город"

# Row 8 - extra key (only the "key" column is populated)
$ws.Range("B8").Value = "system.common.unrequired"

# Match the number-ish column styling used by the existing D/E/F rows
# (center alignment, same cell style used in rows 1/3/4).
$ws.Range("D5:F7").HorizontalAlignment = -4108

# Widen column D so the longer "Street name"/"Neighborhood" labels fit.
$ws.Columns("D").ColumnWidth = 13.166666666666666

# Match the author's final selection.
[void]$ws.Range("B8").Select()
